$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "unit" column wording fix (dollars per person -> dollar per person)
$ws.Range("E4").Value = "dollar per person"

# Row 5: crime_prev_measure / "Short term value to society" updated with
# the new bond-court sourced measure and values.
$ws.Range("S5").Value = "https://loyolaccj.org/blog/cook-bond"
$ws.Range("C5").Value = "Monetary D- bonds measuring the nature of crime, potential danger to public and danger of flee"
$ws.Range("P5").Value = "Central bond court report 2018"
$ws.Range("Q5").Value = "https://cookcountysheriffil.gov/wp-content/uploads/2018/02/Central-Bond-Court-Report.pdf"
$ws.Range("E5").Value = "dollar per person"
$ws.Range("H5").Value = 75000
$ws.Range("J5").Value = 133685
$ws.Range("K5").Value = 75000

# Update the saved selection to match the author's final cursor position.
$ws.Range("R5").Select()
